$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values that look numeric stay as text, matching the source data
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "71.954.91"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "4.024.45"
$ws.Range("E3").Value = "  -0.71%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "534.40"
$ws.Range("E5").Value = "  +1.15%  "

$ws.Range("D6").Value = "150.17"
$ws.Range("E6").Value = "  -0.67%  "

$ws.Range("D7").Value = "4.020.01"
$ws.Range("E7").Value = "  -0.54%  "

$ws.Range("D8").Value = "0.696"
$ws.Range("E8").Value = "  -1.98%  "

$ws.Range("E9").Value = "  -0.07%  "

$ws.Range("E10").Value = "  -1.91%  "

$ws.Range("D11").Value = "0.172"
$ws.Range("E11").Value = "  -2.90%  "

$ws.Range("D12").Value = "54.08"
$ws.Range("E12").Value = "  +6.67%  "

$ws.Range("E13").Value = "  -2.74%  "

$ws.Range("D14").Value = "10.80"
$ws.Range("E14").Value = "  -2.39%  "

$ws.Range("D15").Value = "4.659.41"
$ws.Range("E15").Value = "  -0.88%  "

$ws.Range("D16").Value = "4.018.10"
$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("D17").Value = "14.15"
$ws.Range("E17").Value = "  -1.88%  "

$ws.Range("E18").Value = "  -1.56%  "

$ws.Range("E19").Value = "  -3.49%  "

$ws.Range("E20").Value = "  -1.22%  "

$ws.Range("D21").Value = "71.920.39"
$ws.Range("E21").Value = "  -0.41%  "

$ws.Range("D22").Value = "432.64"
$ws.Range("E22").Value = "  -1.71%  "

$ws.Range("D23").Value = "98.29"
$ws.Range("E23").Value = "  -2.73%  "

$ws.Range("D24").Value = "3.61"
$ws.Range("E24").Value = "  -1.27%  "

$ws.Range("D25").Value = "14.76"
$ws.Range("E25").Value = "  -2.35%  "

$ws.Range("D26").Value = "4.22"
$ws.Range("E26").Value = "  -0.80%  "

$ws.Range("D27").Value = "4.40"
$ws.Range("E27").Value = "  +30.26%  "

$ws.Range("D28").Value = "11.40"
$ws.Range("E28").Value = "  -1.79%  "

$ws.Range("E29").Value = "  -2.75%  "

$ws.Range("D30").Value = "5.95"
$ws.Range("E30").Value = "  +1.89%  "

$ws.Range("D31").Value = "37.01"
$ws.Range("E31").Value = "  -1.47%  "

$ws.Range("E32").Value = "  +21.54%  "

$ws.Range("E33").Value = "  +2.20%  "

$ws.Range("D34").Value = "50.21"
$ws.Range("E34").Value = "  +17.17%  "

$ws.Range("D35").Value = "13.58"
$ws.Range("E35").Value = "  -1.19%  "

$ws.Range("D36").Value = "678.08"
$ws.Range("E36").Value = "  +0.18%  "

$ws.Range("D37").Value = "67.66"
$ws.Range("E37").Value = "  +1.41%  "

$ws.Range("D38").Value = "0.457"
$ws.Range("E38").Value = "  +3.37%  "

$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "3.48"
$ws.Range("E39").Value = "  +9.71%  "

$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0825"
$ws.Range("E40").Value = "  -5.41%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.149"
$ws.Range("E41").Value = "  -5.63%  "

$ws.Range("D42").Value = "3.38"
$ws.Range("E42").Value = "  -2.77%  "

$ws.Range("E43").Value = "  +0.16%  "

$ws.Range("D44").Value = "11.13"
$ws.Range("E44").Value = "  +16.45%  "

$ws.Range("D45").Value = "0.0494"
$ws.Range("E45").Value = "  -2.84%  "

$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("E47").Value = "  -3.03%  "

$ws.Range("D48").Value = "2.65"
$ws.Range("E48").Value = "  -5.56%  "

$ws.Range("E49").Value = "  -0.13%  "

$ws.Range("D50").Value = "3.09"
$ws.Range("E50").Value = "  -0.21%  "

$ws.Range("D51").Value = "2.853.50"
$ws.Range("E51").Value = "  +8.95%  "
